$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row text (capitalization / renamed configuration labels)
$ws.Range("H1").Value = "Population genomics"
$ws.Range("F1").Value = "Biocode"
$ws.Range("C1").Value = "GeOMe"
$ws.Range("D1").Value = "Single sheet generic"
$ws.Range("E1").Value = "Multi sheet generic"
$ws.Range("G1").Value = "Individual-based"

# Mark decimalLatitude / decimalLongitude rows as also required for the
# Biocode (F) and Population genomics (H) configurations
$ws.Range("F6").Value = "required"
$ws.Range("H6").Value = "required"
$ws.Range("F7").Value = "required"
$ws.Range("H7").Value = "required"

# Update the active selection to match the author's final cursor position
$ws.Range("G8").Select()

$wb.Save()
